$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 23, shifting existing rows 23-28 down to 24-29
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with data (mirrors the rest of the block)
$ws.Cells.Item(23, 1).Value = 7
$ws.Cells.Item(23, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(23, 3).Value = "Ñuble"
$ws.Cells.Item(23, 4).Value = 44522
$ws.Cells.Item(23, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(23, 5).Value = 16
$ws.Cells.Item(23, 6).Value = 100112026
$ws.Cells.Item(23, 7).Value = "Haba"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 100
$ws.Cells.Item(23, 11).Value = 6000
$ws.Cells.Item(23, 12).Value = 7000
$ws.Cells.Item(23, 13).Value = 6500
$ws.Cells.Item(23, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(23, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(23, 16).Value = 260
$ws.Cells.Item(23, 17).Value = 25
$ws.Cells.Item(23, 18).Value = "Hortaliza"
